$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 5, pushing the existing "Mean success rate ="
# summary row (D5/E5) down to row 6, leaving row 5 empty.
$ws.Rows.Item(5).Insert()

# --- Row 1 headers: update two existing headers, add new ones for F:J ---
$ws.Range("D1").Value = "Number of attempted auths"
$ws.Range("E1").Value = "Success rate (across authentication)"
$ws.Range("F1").Value = "Success rate (across attempt)"
$ws.Range("G1").Value = "Max Time"
$ws.Range("H1").Value = "Min Time"
$ws.Range("I1").Value = "Mean Time (sec)"
$ws.Range("J1").Value = "SD"

# --- Row 2 (Cube) new measurement columns ---
$ws.Range("F2").Value = 15
$ws.Range("G2").Value = 19.06
$ws.Range("H2").Value = 2.03
$ws.Range("I2").Value = 8.09
$ws.Range("J2").Value = 3.67

# --- Row 3 (Credit Card) new measurement columns ---
$ws.Range("F3").Value = 53.66
$ws.Range("G3").Value = 15.51
$ws.Range("H3").Value = 2.65
$ws.Range("I3").Value = 5.16
$ws.Range("J3").Value = 2.68

# --- Row 4 (Pendant) new measurement columns ---
$ws.Range("F4").Value = 23.4
$ws.Range("G4").Value = 12.29
$ws.Range("H4").Value = 2.29
$ws.Range("I4").Value = 5.47
$ws.Range("J4").Value = 3.49

# --- Column widths for the new/resized columns ---
# (iron_native's ColumnWidth setter only has 1/6-character resolution, so
# these are the closest achievable approximations of the authored widths
# 23.44140625 / 30.33203125 / 24.6640625 / 14.109375.)
$ws.Columns.Item(4).ColumnWidth = 22.666666666666668
$ws.Columns.Item(5).ColumnWidth = 29.5
$ws.Columns.Item(6).ColumnWidth = 23.833333333333332
$ws.Columns.Item(9).ColumnWidth = 13.333333333333334

# --- Restore the selection to match the saved workbook state ---
$ws.Range("E10").Select()
